# Add two new worksheets "non-public" and "with-public" at the end of the
# workbook, replicating the layout/style of the existing regression-output
# sheets ("2016-2019" / "2020-2022") and filling them with the data below.

$wb = $excel.ActiveWorkbook

# --- Reference sheet used as the style/formatting template -----------------
$template = $wb.Worksheets.Item(1)

# --- Row data for "non-public" (becomes sheet3, after "2020-2022") --------
$rows3 = @(
    @($null, "betas", "name", "std", "tval", "pval"),
    @([double]"0", [double]"0.0637088731313248", "end_price_pers", [double]"0.02238315851785786", [double]"2.846286107498913", [double]"0.004423243349374988"),
    @([double]"1", [double]"0.7593975990057419", "pop", [double]"0.1025429604858371", [double]"7.405653156567756", [double]"1.305065219678907e-13"),
    @([double]"2", [double]"0.01251956201842215", "light", [double]"0.002784357840487596", [double]"4.496391173711252", [double]"6.9116529259573e-06"),
    @([double]"3", [double]"-0.01077524367234719", "pm25", [double]"0.001653920264806647", [double]"-6.514971671628246", [double]"7.27033777922428e-11"),
    @([double]"4", [double]"-0.03590023213760504", "elect_store", [double]"0.007784953465750745", [double]"-4.611489624895656", [double]"3.997936150470529e-06"),
    @([double]"5", [double]"0.08042991992061747", "kind", [double]"0.01170855419272161", [double]"6.869329773492884", [double]"6.450423702714271e-12"),
    @([double]"6", [double]"0.0009188494730057676", "age", [double]"0.003505818216727482", [double]"0.2620927316258487", [double]"0.7932499502947528"),
    @([double]"7", [double]"-0.09957865546870498", "hotel_num", [double]"0.01872490744341301", [double]"-5.317978514427019", [double]"1.049264555773517e-07"),
    @([double]"8", [double]"0.0897335773351611", "mall", [double]"0.0100767496612824", [double]"8.905012067526274", [double]"5.338024012456652e-19"),
    @([double]"9", [double]"0.0681412197591295", "museum_num", [double]"0.03207401752967017", [double]"2.124499049615323", [double]"0.03362843543888307"),
    @([double]"10", [double]"0.06201438711410339", "old", [double]"0.03361405766034116", [double]"1.844894411163867", [double]"0.06505289983273083"),
    @([double]"11", [double]"0.01886967576864667", "ktv", [double]"0.006952403994776279", [double]"2.714122450712658", [double]"0.006645159837183084"),
    @([double]"12", [double]"0.05759735098984133", "mid", [double]"0.02495306768315119", [double]"2.308227257714377", [double]"0.0209864998688624"),
    @([double]"13", [double]"0.04044446848911891", "primary", [double]"0.02384813478090232", [double]"1.695917473659492", [double]"0.08990151000025703"),
    @([double]"14", [double]"-0.0007188371302023358", "west_food", [double]"0.009316348107866395", [double]"-0.07715868083496957", [double]"0.9384973115290606"),
    @([double]"15", [double]"0.07676307722657912", "super", [double]"0.01476243462578291", [double]"5.199892780050705", [double]"1.994035203801389e-07"),
    @([double]"16", [double]"-0.1971501125630096", "green_ratio", [double]"0.2171865486487065", [double]"-0.9077455016880198", [double]"0.3640126977409962"),
    @([double]"17", [double]"-0.001735985817913426", "number_building", [double]"0.001487015025690822", [double]"-1.167429910203456", [double]"0.2430367980400935"),
    @([double]"18", [double]"0.05453573902918429", "tihu", [double]"0.1014216886986797", [double]"0.5377127883485363", [double]"0.5907753490410359"),
    @([double]"19", [double]"0.20222518177387", "sub", [double]"0.05122367266021394", [double]"3.947885250542388", [double]"7.884457746141509e-05"),
    @([double]"20", [double]"0.001203659398459148", "floor_ratio", [double]"0.005907814001665293", [double]"0.2037402325326866", [double]"0.8385565005398621"),
    @([double]"21", [double]"1.793027881369573e-05", "residence", [double]"4.237864787995816e-05", [double]"0.4230970007463445", [double]"0.6722244905542698"),
    @([double]"22", [double]"0.1050286737796948", "park", [double]"0.01611663819562338", [double]"6.516785480002669", [double]"7.183003569083262e-11"),
    @([double]"23", [double]"0.3874334265550495", "W_kou", [double]"0.01364562403362328", [double]"28.39250338426447", [double]"2.502484616927769e-177")
  )

# --- Row data for "with-public" (becomes sheet4, after "non-public") ------
$rows4 = @(
    @($null, "betas", "name", "std", "tval", "pval"),
    @([double]"0", [double]"0.04377605336263387", "end_price_pers", [double]"0.00564526274949352", [double]"7.754475797705139", [double]"8.870913996172622e-15"),
    @([double]"1", [double]"0.959082868500801", "pop", [double]"0.04749848506487064", [double]"20.19186227078489", [double]"1.154305096726273e-90"),
    @([double]"2", [double]"0.01073859511211014", "light", [double]"0.0008143977526359257", [double]"13.18593411800683", [double]"1.057361136058038e-39"),
    @([double]"3", [double]"-0.005749974025439229", "pm25", [double]"0.0004883022089631953", [double]"-11.77544135556557", [double]"5.224294340831422e-32"),
    @([double]"4", [double]"-0.06721149416600077", "elect_store", [double]"0.005584747257481041", [double]"-12.03483184954658", [double]"2.331084923773761e-33"),
    @([double]"5", [double]"0.09619071195437488", "kind", [double]"0.004008706111860803", [double]"23.99545121798018", [double]"3.10210932730317e-127"),
    @([double]"6", [double]"0.005449340192978054", "age", [double]"0.001094841438308703", [double]"4.977287123326393", [double]"6.448163618429017e-07"),
    @([double]"7", [double]"-0.1202579828244073", "hotel_num", [double]"0.009976158486020856", [double]"-12.05453812636592", [double]"1.835585699084158e-33"),
    @([double]"8", [double]"0.05567348910636791", "mall", [double]"0.00395080415821821", [double]"14.09168535741249", [double]"4.272421689789468e-45"),
    @([double]"9", [double]"0.1303558675702017", "museum_num", [double]"0.0197488292065366", [double]"6.60068838546923", [double]"4.092529716942461e-11"),
    @([double]"10", [double]"0.05613176093359085", "old", [double]"0.01687623348327957", [double]"3.326083452756469", [double]"0.0008807555450991296"),
    @([double]"11", [double]"0.01948603866213898", "ktv", [double]"0.003789171827270449", [double]"5.142558730617359", [double]"2.710216776840407e-07"),
    @([double]"12", [double]"0.0387479538324378", "mid", [double]"0.01126137518572047", [double]"3.440783491661885", [double]"0.0005800324588749407"),
    @([double]"13", [double]"0.1081165504577916", "primary", [double]"0.01083609246706428", [double]"9.977448124072955", [double]"1.913240165307493e-23"),
    @([double]"14", [double]"0.01425896921789705", "west_food", [double]"0.004593741965328372", [double]"3.103998728164911", [double]"0.0019092415168455"),
    @([double]"15", [double]"0.0861367748886903", "super", [double]"0.006196877208041567", [double]"13.90002932072178", [double]"6.33154231032481e-44"),
    @([double]"16", [double]"-0.4016180539608001", "green_ratio", [double]"0.06294278770727448", [double]"-6.380684246598502", [double]"1.762984750860955e-10"),
    @([double]"17", [double]"-0.0002944294891112696", "number_building", [double]"0.0001448205766837595", [double]"-2.03306391849417", [double]"0.04204606572397063"),
    @([double]"18", [double]"-0.04206458033378282", "tihu", [double]"0.01687667229857829", [double]"-2.492468870022817", [double]"0.01268584503278715"),
    @([double]"19", [double]"0.1840730562249867", "sub", [double]"0.02006160458153283", [double]"9.175390506621298", [double]"4.499406579171819e-20"),
    @([double]"20", [double]"5.376666124485034e-05", "floor_ratio", [double]"0.0006539656863377688", [double]"0.08221633392715995", [double]"0.9344746848226322"),
    @([double]"21", [double]"2.616653593235976e-06", "residence", [double]"6.642911969326626e-06", [double]"0.3939015909466009", [double]"0.693653695733659"),
    @([double]"22", [double]"0.08687055406353217", "park", [double]"0.006052228618040066", [double]"14.35348192310423", [double]"1.01319875154469e-46"),
    @([double]"23", [double]"0.3573186051384035", "W_kou", [double]"0.005803869959308016", [double]"61.56557738950545", [double]"0")
  )

function Add-RegressionSheet {
    param(
        [string]$SheetName,
        [object[]]$Rows
    )

    # Add the sheet after the last existing sheet so ordering matches.
    $lastIndex = $wb.Worksheets.Count
    $ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
    $ws.Name = $SheetName

    # Match the outlinePr summaryBelow/summaryRight="1" used on the other
    # sheets in this workbook.
    $ws.Outline.SummaryRow = 1
    $ws.Outline.SummaryColumn = 1

    # Match page margins used by the other sheets in the workbook
    # (left/right 0.75in, top/bottom 1in, header/footer 0.5in).
    $ws.PageSetup.LeftMargin = 54
    $ws.PageSetup.RightMargin = 54
    $ws.PageSetup.TopMargin = 72
    $ws.PageSetup.BottomMargin = 72
    $ws.PageSetup.HeaderMargin = 36
    $ws.PageSetup.FooterMargin = 36

    # Write out all cell values row by row / column by column.
    $r = 1
    foreach ($row in $Rows) {
        for ($c = 0; $c -lt $row.Length; $c++) {
            $val = $row[$c]
            if ($val -ne $null) {
                $ws.Cells.Item($r, $c + 1).Value = $val
            }
        }
        $r++
    }

    # Copy over the same cell formatting used on the template sheet:
    #  - bold, centered, bordered header cells in row 1 (columns B:F)
    #  - bordered/centered index cells in column A (rows 2:25)
    $template.Range("B1:F1").Copy()
    $ws.Range("B1:F1").PasteSpecial(-4122) | Out-Null

    $template.Range("A2:A25").Copy()
    $ws.Range("A2:A25").PasteSpecial(-4122) | Out-Null

    $excel.CutCopyMode = 0

    return $ws
}

$wsNonPublic = Add-RegressionSheet "non-public" $rows3
$wsWithPublic = Add-RegressionSheet "with-public" $rows4

Write-Host "Added sheets: $($wsNonPublic.Name), $($wsWithPublic.Name)"
Write-Host "Total worksheets: $($wb.Worksheets.Count)"
